$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2
$wdReplaceOne = 1
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($old, $new, $replaceMode) {
    $r = $d.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $replaceMode) | Out-Null
}

# "English" appears twice in the doc (hyperlink label + standalone paragraph);
# both translate to "Inglese", so replace all occurrences at once.
Replace-Text "English" "Inglese" $wdReplaceAll

# Language list after the English hyperlink.
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portoghese / Francese / Thailandese / Vietnamita / Spagnolo" $wdReplaceOne

# Brief / summary table.
Replace-Text "Brief" "Riepilogo" $wdReplaceOne
Replace-Text "An email sent upon verification to partners in the target country who have sent the correct documents. It will be sent via customer.io" "Un'email inviata dopo la verifica ai partner nel paese di destinazione che hanno inviato i documenti corretti. Sarà inviata tramite customer.io" $wdReplaceOne
Replace-Text "Target audience" "Pubblico target" $wdReplaceOne
Replace-Text "Invited partners who didn’t submit their documents on time" "Partner invitati che non hanno presentato i loro documenti in tempo" $wdReplaceOne

# Subject line.
Replace-Text "Subject line" "Oggetto" $wdReplaceOne
Replace-Text " — one step closer!" " — un passo più vicino!" $wdReplaceOne

# Heading + greeting.
Replace-Text "Your documents have been verified!" "I tuoi documenti sono stati verificati!" $wdReplaceOne
Replace-Text "Hi " "Ciao " $wdReplaceOne

# Body paragraphs.
Replace-Text "We’ve reviewed the documents you’ve sent us for the " "Abbiamo esaminato i documenti che ci hai inviato per l'evento " $wdReplaceOne
Replace-Text " and all of them have been verified! " " e tutti sono stati verificati! " $wdReplaceOne
Replace-Text "We’ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly." "Presto invieremo ulteriori dettagli sull'evento, tra cui l'agenda e le modalità di viaggio. Pertanto, assicurati di controllare regolarmente la casella di posta elettronica." $wdReplaceOne

# Contact-us paragraph (live chat / WhatsApp).
Replace-Text "If you have any questions, please contact us via " "Se hai domande, non esitare a contattarci tramite " $wdReplaceOne
Replace-Text "live chat" "chat live" $wdReplaceOne
# First " or " (between live chat and WhatsApp links) -> " o "
Replace-Text " or " " o " $wdReplaceOne

# Country-manager paragraph.
Replace-Text "If you have any questions, please contact your country manager, " "Se hai domande, contatta il tuo country manager, " $wdReplaceOne
Replace-Text ", at " ", all'indirizzo " $wdReplaceOne
# Second " or " (between email address and WhatsApp number) -> " o al numero "
Replace-Text " or " " o al numero " $wdReplaceOne

# Comment text.
$c = $d.Comments.Item(1)
$c.Range.Text = "scegli uno dei due"
